# 0.2 noozle miniatures profile added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custos")

# Update the input parameters for the new 0.2 nozzle miniatures profile
$ws.Range("C6").Value = 2     # Tempo de impressao por unidade - horas
$ws.Range("C7").Value = 37    # Tempo de impressao por unidade - minutos
$ws.Range("C8").Value = 4     # Peso por unidade
$ws.Range("C25").Value = 20   # Percentual de falhas de impressao

# Reflect the selection change noted in the saved view state
$ws.Range("C30").Select()

$wb.Save()
